# "Generate Report for Handback" — the handback files are now in sync with
# en-US, so the per-language Status is updated, the Latest Handback
# DateTime is refreshed, and the stale "version not latest" Error Detail
# is cleared out on both the zh-cn and de-de report rows (and the rollup
# Overview row that mirrors the per-language status).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("L2").Value = "2017-02-09 13:48:32"
$wsZhCn.Range("R2").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("L2").Value = "2017-02-09 13:48:55"
$wsDeDe.Range("R2").Value = ""

# --- Column width adjustments to fit the new, longer Status text and the
#     now-empty Error Detail column ---
$wideWidth = 29.166666666666668   # -> stored width ~29.98 (Status columns)
$narrowWidth = 12.833333333333334 # -> stored width ~13.75 (Error Detail columns)

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(18).ColumnWidth = $narrowWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(18).ColumnWidth = $narrowWidth
